$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update revised "Positive New" (C column) values for rows with corrected case counts.
# B (cumulative total) and D (7-day average) columns are formulas that will recalculate automatically.
$ws.Range("C244").Value = 225
$ws.Range("C267").Value = 1120
$ws.Range("C280").Value = 1837
$ws.Range("C285").Value = 1327
$ws.Range("C288").Value = 2821
$ws.Range("C289").Value = 2659
$ws.Range("C294").Value = 3525
$ws.Range("C300").Value = 1193
$ws.Range("C301").Value = 3590
$ws.Range("C302").Value = 3788
$ws.Range("C308").Value = 5503
$ws.Range("C309").Value = 5858
$ws.Range("C310").Value = 6073
$ws.Range("C311").Value = 5760
$ws.Range("C315").Value = 6205
$ws.Range("C316").Value = 5391
$ws.Range("C317").Value = 5395
$ws.Range("C318").Value = 5454
$ws.Range("C319").Value = 4936
$ws.Range("C321").Value = 2202
$ws.Range("C322").Value = 6259
$ws.Range("C323").Value = 5998
$ws.Range("C324").Value = 5634
$ws.Range("C325").Value = 1564
$ws.Range("C326").Value = 5655
$ws.Range("C327").Value = 3634
$ws.Range("C328").Value = 2305
$ws.Range("C329").Value = 6449
$ws.Range("C330").Value = 5933
$ws.Range("C331").Value = 4464
$ws.Range("C332").Value = 2618
$ws.Range("C333").Value = 489
$ws.Range("C334").Value = 4137

# Append new daily rows (335-340) through Jan 1, 2021, continuing the table
# by cloning the format/formulas of the last existing row and inserting below it.
$lastRow = 334
$newRowsData = @(
    @{Date = 44192; PositiveNew = 2637},
    @{Date = 44193; PositiveNew = 8137},
    @{Date = 44194; PositiveNew = 6405},
    @{Date = 44195; PositiveNew = 3919},
    @{Date = 44196; PositiveNew = 1083},
    @{Date = 44197; PositiveNew = 96}
)
foreach ($nr in $newRowsData) {
    $newRowNum = $lastRow + 1
    $ws.Range("A" + $lastRow + ":D" + $lastRow).Copy() | Out-Null
    $ws.Range("A" + $newRowNum + ":D" + $newRowNum).Insert(-4121) | Out-Null
    $ws.Range("A" + $newRowNum).Value = $nr.Date
    $ws.Range("C" + $newRowNum).Value = $nr.PositiveNew
    $lastRow = $newRowNum
}
$excel.CutCopyMode = $false

# Update the active selection to reflect where the editor last left off in the sheet.
$ws.Range("J336").Select() | Out-Null
